# Change _id field generator, added path for MapRDB

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Sheet1 data ---
# H2 / I2 (map.parallelism / shuffle.parallelism) for Maprdb-aggregation row -> "-"
$ws1.Range("H2").Value = "-"
$ws1.Range("I2").Value = "-"

# F3: Maprdb-scan data size uservisits count changed from 5000000000 to 2000000000
$ws1.Range("F3").Value = "hibench.join.bigdata.uservisits                 2000000000" + [char]10 + "hibench.join.bigdata.pages                      120000000"

# Row 3 shrinks (wrapped text now fits in fewer lines after the edit)
$ws1.Rows.Item(3).RowHeight = 25.7

# --- Add Sheet2 (new empty worksheet), placed after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("D21").Select()

# Keep Sheet1 as the active/selected sheet & restore its selection
$ws1.Activate()
$ws1.Range("G25").Select()
